# Regenerate the s_vals data to filter save games: update B2:E13 with the
# newly computed per-game values, and recompute the "sum" column G as the
# row total of B+C+D+E (column F "Win" is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 3.286832544864788;  C = 1.655778082260271;   D = 22.3905356188092;   E = 10.19245300693656 }
    3  = @{ B = 1.455362044514542;  C = 1.655778082260271;   D = 0.7527432677738641; E = 0.4942365360607697 }
    4  = @{ B = 0.1190320826869504; C = 0.306821227259698;    D = 3.537761648806719;  E = 0.4942365360607697 }
    5  = @{ B = 0.1190320826869504; C = 0.306821227259698;    D = 0.7527432677738641; E = 0.4942365360607697 }
    6  = @{ B = 0.6606524410359556; C = 0.04071648406533734;  D = 22.3905356188092;   E = 0.4942365360607697 }
    7  = @{ B = 3.286832544864788;  C = 1.655778082260271;    D = 3.537761648806719;  E = 0.4942365360607697 }
    8  = @{ B = 3.286832544864788;  C = 1.655778082260271;    D = 22.3905356188092;   E = 0.4942365360607697 }
    9  = @{ B = 0.2917716402565462; C = 1.655778082260271;    D = 3.537761648806719;  E = 0.4942365360607697 }
    10 = @{ B = 1.455362044514542;  C = 1.655778082260271;    D = 0.7527432677738641; E = 1133.036916526867 }
    11 = @{ B = 0.6606524410359556; C = 0.306821227259698;    D = 0.7527432677738641; E = 0.4942365360607697 }
    12 = @{ B = 1.455362044514542;  C = 1.655778082260271;    D = 3.537761648806719;  E = 0.4942365360607697 }
    13 = @{ B = 1.455362044514542;  C = 1.655778082260271;    D = 3.537761648806719;  E = 0.4942365360607697 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 7).Value = ($row.B + $row.C + $row.D + $row.E)
}
